$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining worker row (row 16) with the new worker's data
$ws.Range("C16").Value = "1056779330"
$ws.Range("D16").Value = "HARRINSON ALIRIO RUIZ BLANDON"

# Delete row 17 (the old second worker row), shifting rows below up
$ws.Rows("17").Delete()

# Update "Cant. Trabajadores" (Number of workers) from 2 to 1
$ws.Range("C13").Value = 1

# Update "Valor Mora" total from 66250 to 33125
$ws.Range("E11").Value = 33125
